$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("class_schedule")

# Row 26: "Parallelism / Distributed Computing" week
# - drop "Part 1" from topic, since now spread over two weeks
$ws.Range("B26").Value = "- Parallelism `n- Distributed Computing"
# - renumber the two exercise links
$ws.Range("D26").Value = "- ``Link 1 <exercises/Exercise_dask.ipynb>```_`n- ``Link 2 <exercises/Exercise_dask_realdata.ipynb>```_"

# Row 27: "Distributed Computing, Azure" week - add in-class exercise link
$ws.Range("D27").Value = "``Link <exercises/Exercise_azure_arcos.ipynb>```_"

# Row 28: "Distributed Computing, Azure 2" week - opioid project due notice now
# becomes a bulleted list, with two new reading assignments added
$ws.Range("C28").Value = "- **OPIOID PROJECT DUE (Extensions available upon request to assigned date for PDS final)**`n- ``Azure Storage <cloud_azurestorage.ipynb>```_`n- ``More Azure Concepts <cloud_more_concepts.ipynb>```_"

# Update the last active selection on the sheet to C29 (matches saved view state)
$ws.Range("C29").Select()
